$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp cell (A1): 10:05 -> 10:35
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 10:35"

# --- Reorder country: move "Consejo Danes para los Refugiados" row up to where
# "Republica de Yibuti" currently is (row 88), shifting the rows below it down
# by one (Yibuti -> 89, Islandia -> 90, Estonia -> 91), and refresh the
# latest case numbers for the affected rows.

$ws.Range("A88").Value = "Consejo Danes para los Refugiados"
$ws.Range("B88").Value = 1835
$ws.Range("C88").Value = 104
$ws.Range("D88").Value = 303
$ws.Range("E88").Value = 1471
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 61

$ws.Range("A89").Value = "Republica de Yibuti"
$ws.Range("B89").Value = 1828
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 1052
$ws.Range("E89").Value = 767
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 9

$ws.Range("A90").Value = "Islandia"
$ws.Range("B90").Value = 1803
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 1789
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 10

$ws.Range("A91").Value = "Estonia"
$ws.Range("B91").Value = 1800
$ws.Range("C91").Value = 6
$ws.Range("D91").Value = 963
$ws.Range("E91").Value = 773
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 64

# --- Updated case counts for other countries (Polonia, Israel, Filipinas, Montenegro)

# Polonia
$ws.Range("B34").Value = 19983
$ws.Range("C34").Value = 244
$ws.Range("D34").Value = 8452
$ws.Range("E34").Value = 10566
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 3
$ws.Range("H34").Value = 965

# Israel
$ws.Range("B41").Value = 16670
$ws.Range("C41").Value = 3
$ws.Range("D41").Value = 13617
$ws.Range("E41").Value = 2774
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 279

# Filipinas
$ws.Range("B46").Value = 13434
$ws.Range("C46").Value = 213
$ws.Range("D46").Value = 3000
$ws.Range("E46").Value = 9588
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 4
$ws.Range("H46").Value = 846

# Montenegro
$ws.Range("D144").Value = 314
$ws.Range("E144").Value = 1
